$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.355.02"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "3.006.38"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'554.91"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "'152.70"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.79%  "
$ws.Range("D9").Value = "3.010.12"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -5.87%  "
$ws.Range("D12").Value = "'0.368"
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("D13").Value = "3.528.74"
$ws.Range("E13").Value = "  -4.40%  "
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("D15").Value = "62.488.01"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'23.81"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "3.007.06"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").Value = "'394.59"
$ws.Range("E19").Value = "  -3.92%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'12.02"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'65.16"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").Value = "'0.469"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("E26").Value = "  -6.60%  "
$ws.Range("D27").Value = "0.0₃0974"
$ws.Range("E27").Value = "  -4.92%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.54"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("D32").Value = "'20.57"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").Value = "'160.54"
$ws.Range("E33").Value = "  +4.79%  "
$ws.Range("D34").Value = "'4.71"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("E36").Value = "  -2.89%  "
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  -5.25%  "
$ws.Range("D39").Value = "2.458.24"
$ws.Range("E39").Value = "  -10.48%  "
$ws.Range("D40").Value = "'3.94"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "'22.60"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").Value = "'37.54"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "'0.663"
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("D44").Value = "'0.0599"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").Value = "'0.0249"
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  -9.47%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'19.90"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0956"
$ws.Range("E49").Value = "  -2.10%  "
$ws.Range("D50").Value = "'10.49"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "'264.68"
$ws.Range("E51").Value = "  -6.97%  "
